$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: give the three brand-new rows (5,6,7) the same look (borders,
# fill, font, alignment) as the existing data rows, by copying row 2's
# formatting down. Doing this before any value/number-format changes keeps
# every data row on a consistent style baseline. ---
$ws.Range("A2:K2").Copy() | Out-Null
$ws.Range("A5:K5").PasteSpecial(-4122) | Out-Null
$ws.Range("A6:K6").PasteSpecial(-4122) | Out-Null
$ws.Range("A7:K7").PasteSpecial(-4122) | Out-Null

# --- Step 2: the "Date" column (A) stores dates as plain text
# (e.g. "31-JAN-26"), not real date serials. Force text format on the whole
# column range first so typing a date-like string below doesn't get
# auto-converted into a date value/format. ---
$ws.Range("A2:A7").NumberFormat = "@"

# ---- Row 2 : Air Arabia Egypt E5-592 ----
$ws.Range("A2").Value = "28-JAN-26"
$ws.Range("B2").Value = "SM-436"
$ws.Range("C2").Value = "Air Arabia Egypt E5-592"
$ws.Range("D2").Value = 345
$ws.Range("E2").Value = 602
$ws.Range("F2").Value = -257
$ws.Range("G2").Value = 30
$ws.Range("H2").Value = 30
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = "LOW THREAT"
$ws.Range("K2").Value = "SAR"

# ---- Row 3 : Nile Air NP-106 ----
$ws.Range("A3").Value = "28-JAN-26"
$ws.Range("B3").Value = "SM-436"
$ws.Range("C3").Value = "Nile Air NP-106"
$ws.Range("D3").Value = 350
$ws.Range("E3").Value = 602
$ws.Range("F3").Value = -252
$ws.Range("G3").Value = 30
$ws.Range("H3").Value = 30
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = "LOW THREAT"
$ws.Range("K3").Value = "SAR"

# ---- Row 4 : Nile Air NP-116 (previously EgyptAir MS-634 / 07-FEB-26) ----
$ws.Range("A4").Value = "31-JAN-26"
$ws.Range("B4").Value = "SM-436"
$ws.Range("C4").Value = "Nile Air NP-116"
$ws.Range("D4").Value = 345
$ws.Range("E4").Value = 602
$ws.Range("F4").Value = -257
$ws.Range("G4").Value = 30
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = "LOW THREAT"
$ws.Range("K4").Value = "SAR"

# ---- Row 5 (new) : Nesma Airlines NE-141 ----
$ws.Range("A5").Value = "31-JAN-26"
$ws.Range("B5").Value = "SM-436"
$ws.Range("C5").Value = "Nesma Airlines NE-141"
$ws.Range("D5").Value = 350
$ws.Range("E5").Value = 602
$ws.Range("F5").Value = -252
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = 30
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = "LOW THREAT"
$ws.Range("K5").Value = "SAR"

# ---- Row 6 (new) : Nesma Airlines NE-141 ----
$ws.Range("A6").Value = "14-FEB-26"
$ws.Range("B6").Value = "SM-436"
$ws.Range("C6").Value = "Nesma Airlines NE-141"
$ws.Range("D6").Value = 400
$ws.Range("E6").Value = 663
$ws.Range("F6").Value = -263
$ws.Range("G6").Value = 30
$ws.Range("H6").Value = 30
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = "LOW THREAT"
$ws.Range("K6").Value = "SAR"

# ---- Row 7 (new) : Nile Air NP-116 ----
$ws.Range("A7").Value = "14-FEB-26"
$ws.Range("B7").Value = "SM-436"
$ws.Range("C7").Value = "Nile Air NP-116"
$ws.Range("D7").Value = 401
$ws.Range("E7").Value = 663
$ws.Range("F7").Value = -262
$ws.Range("G7").Value = 30
$ws.Range("H7").Value = 30
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = "LOW THREAT"
$ws.Range("K7").Value = "SAR"
